# financas.xlsx: adiciona duas novas transações (Despesa / Receita) à planilha
# para que o saldo possa ser calculado como a diferença entre receitas e despesas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Linha 2 - Despesa (cartão de crédito), com data/hora 2024-06-28
$ws.Range("A2").Value = 45471
$ws.Range("B2").Value = "Despesa"
$ws.Range("C2").Value = "Cartão de Crédito"
$ws.Range("D2").Value = 1400
$ws.Range("E2").Value = "cartao sofisa"
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Linha 3 - Receita (salário), com data 2024-06-28
$ws.Range("A3").Value = 45471
$ws.Range("B3").Value = "Receita"
$ws.Range("C3").Value = "Salário"
$ws.Range("D3").Value = 4000
$ws.Range("E3").Value = ""
$ws.Range("A3").NumberFormat = "yyyy-mm-dd"
$ws.Range("A3").NumberFormat = "YYYY-MM-DD"
